# Auto-generated Excel COM-interop edit script
# Implements the scraping-refresh update for Linea 141 (horarios-141-2026-01-15.xlsx)
# covering sheets "LP1912" (1), "LP1912-215" (2) and "6203-6173" (3).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "LP1912"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Header metadata: timestamp + row count
$ws1.Range("A2").Value = "Última actualización: 15:02:32"
$ws1.Range("A3").Value = "Total filas: 289"

# Data rows that changed (reordered / new rows from the refreshed scrape).
# Columns: Hora_Scrap(A) Hora_Llegada(B) Linea(C) Minutos(D) Parada(E)
$sheet1Data = @(
    @(86, "08:55:19", "08:56", "16_SANTA ANA", 1, "LP1912"),
    @(87, "08:55:19", "08:56", "10_OLMOS", 1, "LP1912"),
    @(140, "10:52:48", "10:52", "10_OLMOS", 0, "LP1912"),
    @(141, "10:12:35", "10:52", "15_ABASTO", 40, "LP1912"),
    @(177, "10:12:35", "11:59", "225_GOMEZ", 107, "LP1912"),
    @(178, "11:59:06", "11:59", "16_SANTA ANA", 0, "LP1912"),
    @(186, "11:17:08", "12:07", "16_P MOR-SANTA ANA", 50, "LP1912"),
    @(187, "11:17:08", "12:07", "14_ABASTO", 50, "LP1912"),
    @(188, "11:59:06", "12:07", "84_COLONIA URQUIZA-ESC 49", 8, "LP1912"),
    @(202, "11:45:01", "12:35", "23_HERNANDEZ", 50, "LP1912"),
    @(203, "11:45:01", "12:35", "11_ETCHEVERRY", 50, "LP1912"),
    @(240, "14:00:52", "14:00", "14_ABASTO", 0, "LP1912"),
    @(241, "14:00:52", "14:00", "16_SANTA ANA", 0, "LP1912"),
    @(260, "15:02:32", "15:03", "15_ABASTO", 1, "LP1912"),
    @(262, "14:00:52", "15:04", "10_OLMOS", 64, "LP1912"),
    @(263, "13:23:09", "15:05", "10_OLMOS", 102, "LP1912"),
    @(264, "14:44:25", "15:06", "16_SANTA ANA", 22, "LP1912"),
    @(265, "14:00:52", "15:10", "17_ROMERO", 70, "LP1912"),
    @(266, "14:00:52", "15:13", "11_ETCHEVERRY", 73, "LP1912"),
    @(267, "13:23:09", "15:14", "11_ETCHEVERRY", 111, "LP1912"),
    @(268, "14:44:25", "15:16", "16_SANTA ANA", 32, "LP1912"),
    @(269, "14:00:52", "15:20", "15_ABASTO", 80, "LP1912"),
    @(270, "13:23:09", "15:21", "26_HERNANDEZ", 118, "LP1912"),
    @(271, "14:00:52", "15:25", "26_HERNANDEZ", 85, "LP1912"),
    @(272, "15:02:32", "15:26", "16_SANTA ANA", 24, "LP1912"),
    @(273, "14:00:52", "15:32", "84_COLONIA URQUIZA-ESC 49", 92, "LP1912"),
    @(274, "15:02:32", "15:34", "23_HERNANDEZ", 32, "LP1912"),
    @(275, "14:00:52", "15:35", "23_HERNANDEZ", 95, "LP1912"),
    @(276, "14:00:52", "15:36", "10_OLMOS", 96, "LP1912"),
    @(277, "14:44:25", "15:37", "10_OLMOS", 53, "LP1912"),
    @(278, "14:00:52", "15:38", "215A_EL PATO", 98, "LP1912"),
    @(279, "14:44:25", "15:45", "14_ABASTO", 61, "LP1912"),
    @(280, "14:00:52", "15:46", "16_P MOR-167 Y 521", 106, "LP1912"),
    @(281, "14:44:25", "15:47", "23_HERNANDEZ", 63, "LP1912"),
    @(282, "14:00:52", "15:53", "11_ETCHEVERRY", 113, "LP1912"),
    @(283, "14:44:25", "15:56", "17_ROMERO", 72, "LP1912"),
    @(284, "14:00:52", "15:56", "27_EL RETIRO", 116, "LP1912"),
    @(285, "14:44:25", "16:01", "10_OLMOS", 77, "LP1912"),
    @(286, "15:02:32", "16:01", "27_EL RETIRO", 59, "LP1912"),
    @(287, "14:44:25", "16:02", "27_EL RETIRO", 78, "LP1912"),
    @(288, "15:02:32", "16:08", "14_ABASTO", 66, "LP1912"),
    @(289, "14:44:25", "16:15", "225_C ROCA-H SUR", 91, "LP1912"),
    @(290, "14:44:25", "16:20", "215C_EL PATO", 96, "LP1912"),
    @(291, "14:44:25", "16:21", "26_HERNANDEZ", 97, "LP1912"),
    @(292, "14:44:25", "16:42", "16_P MOR-SANTA ANA", 118, "LP1912"),
    @(293, "14:44:25", "16:43", "225_GOMEZ", 119, "LP1912"),
    @(294, "15:02:32", "16:56", "17_179 Y 38", 114, "LP1912")
)

foreach ($row in $sheet1Data) {
    $r = $row[0]
    $ws1.Cells.Item($r, 1).Value = $row[1]
    $ws1.Cells.Item($r, 2).Value = $row[2]
    $ws1.Cells.Item($r, 3).Value = $row[3]
    $ws1.Cells.Item($r, 4).Value = $row[4]
    $ws1.Cells.Item($r, 5).Value = $row[5]
}

# ---------------------------------------------------------------------------
# Sheet 2: "LP1912-215"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A2").Value = "Última actualización: 15:02:32"

# ---------------------------------------------------------------------------
# Sheet 3: "6203-6173"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("A2").Value = "Última actualización: 15:02:32"
$ws3.Range("A3").Value = "Total filas: 40"

# New row 45
$ws3.Cells.Item(45, 1).Value = "15:02:32"
$ws3.Cells.Item(45, 2).Value = "16:53"
$ws3.Cells.Item(45, 3).Value = "215B_LP-P MOR-40 Y 115"
$ws3.Cells.Item(45, 4).Value = 111
$ws3.Cells.Item(45, 5).Value = "L6173"
